$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(4, 5, 6, 9, 18, 22, 32, 38, 39, 40, 41, 42, 43)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = "T"
}
